$d = $word.ActiveDocument

# The page used to end with the final bibliography entry, a blank
# paragraph, a "Ver no Jupiter..." navigation line, and a copyright
# footer line. The site rebuild dropped the footer block, so remove
# the blank paragraph plus the two footer paragraphs that follow the
# last bibliography entry, leaving a single trailing blank paragraph
# before the page break.

$anchorText = "Gestão Empresarial - Estratégias Organizacionais Autor: Bertero, C. O. Editora: ATLAS"

$count = $d.Paragraphs.Count
$anchorIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text.TrimEnd("`r`a")
    if ($text -eq $anchorText) {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -gt 0) {
    $startPara = $d.Paragraphs.Item($anchorIndex + 1)
    $endPara = $d.Paragraphs.Item($anchorIndex + 3)
    $range = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $range.Delete()
}
